{"js": "// Office.js (Word JavaScript API) edit script.\n// 1) Insert a new \"FirstParagraph\"-styled paragraph right after the\n//    \"Research Progress\" heading, describing the two in-progress manuscripts\n//    (GFM paper + acI paper), mirroring the run/break layout from the diff.\n// 2) Split the run-on \"Oral Presentations\" paragraph into two lines (talks),\n//    replacing the single space between the two quoted talk titles with a\n//    line break.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet researchProgressPara = null;\nlet oralPresentationsPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const t = p.text;\n  if (t === \"Research Progress\" && !researchProgressPara) {\n    researchProgressPara = p;\n  }\n  if (t.indexOf('\"Genome-wide Selective Sweeps') === 0 && !oralPresentationsPara) {\n    oralPresentationsPara = p;\n  }\n}\n\nif (!researchProgressPara) {\n  throw new Error(\"Could not find the 'Research Progress' heading paragraph\");\n}\nif (!oralPresentationsPara) {\n  throw new Error(\"Could not find the Oral Presentations talks paragraph\");\n}\n\n// --- 1) New manuscript-status paragraph under \"Research Progress\" ---\nconst newPara = researchProgressPara.insertParagraph(\"\", \"After\");\nnewPara.style = \"FirstParagraph\";\n\nnewPara.insertText(\"I currently have a manuscript submitted\", \"End\");\nnewPara.insertText(\" \", \"End\");\nnewPara.insertText(\" \", \"End\");\nnewPara.insertText(\n  \"on the results of my analysis of 30 genomes from metagenomes(GFMs). \" +\n    \"I mapped the reads from our metagenomic time series of Trout Bog to the GFMs. \" +\n    \"From the mapping I was able to extract coverage and single nucleotide polymorphism(SNP) data through time. \" +\n    \"We found evidence for both genome-wide and gene-specific sweeps depending on the sequence-discrete population investigated.\",\n  \"End\"\n);\nnewPara.insertText(\" \", \"End\");\nnewPara.insertBreak(Word.BreakType.line, \"End\");\nnewPara.insertText(\n  \"I am also currently working on a manuscript concerning the population dynamics of the ubiquitous and abundant freshwater bacterium acI. \" +\n    \"We have 14 single cell genomes(SAGs) from this clade, representing 3 lakes and X tribes. \" +\n    \"Since these genomes come from a single cell and not an assembly of many cells, like the GFMs, we shredded the genomes and mapped them against each other in order to understand how their relationships would hold up in our metagenomic time series. \" +\n    \"We found that SAGs from different tribes do not map well to each other. \" +\n    \"When we mapped reads from the metagenomic time series of Mendota, we found that only SAGs colleted from the same lake were representive of the populations in the lake. \" +\n    \"We also found that members of the same tribe represent ecologically distinct populations in the lake.\",\n  \"End\"\n);\nnewPara.insertText(\" \", \"End\");\nnewPara.insertBreak(Word.BreakType.line, \"End\");\n\n// --- 2) Split the Oral Presentations talk list into two lines ---\noralPresentationsPara.insertText(\n  '\"Genome-wide Selective Sweeps in Natural Bacterial Populations Revealed by Time-series Metagenomics.\" 15th International Symposium on Microbial Ecology. August 24, 2014. Seoul, South Korea',\n  \"Replace\"\n);\noralPresentationsPara.insertBreak(Word.BreakType.line, \"End\");\noralPresentationsPara.insertText(\n  '\"Genome-wide and Gene-specific Selective Sweeps in Freshwater Bacterial Populations Revealed Using Metagenomics.\" 14 Symposium Society for Aquatic Microbial Ecology. August 2015. Uppsala, Sweden',\n  \"End\"\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# 1) Insert a new \"FirstParagraph\"-styled paragraph right after the\n#    \"Research Progress\" heading, describing the two in-progress manuscripts\n#    (GFM paper + acI paper), mirroring the run/break layout from the diff.\n# 2) Split the run-on \"Oral Presentations\" paragraph into two lines (talks),\n#    replacing the single space between the two quoted talk titles with a\n#    line break.\n\n$d = $word.ActiveDocument\n$lineBreak = [char]11\n\n$researchProgress = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Research Progress\" -and $researchProgress -eq $null) {\n        $researchProgress = $p\n    }\n}\n\nif ($researchProgress -eq $null) {\n    throw \"Could not find the 'Research Progress' heading paragraph\"\n}\n\n# --- 1) New manuscript-status paragraph under \"Research Progress\" ---\n$r = $researchProgress.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n$newPara = $researchProgress.Next()\n$newPara.Range.Style = \"FirstParagraph\"\n\n$part1 = \"I currently have a manuscript submitted\"\n$part1 += \"  \"\n$part1 += \"on the results of my analysis of 30 genomes from metagenomes(GFMs). \"\n$part1 += \"I mapped the reads from our metagenomic time series of Trout Bog to the GFMs. \"\n$part1 += \"From the mapping I was able to extract coverage and single nucleotide polymorphism(SNP) data through time. \"\n$part1 += \"We found evidence for both genome-wide and gene-specific sweeps depending on the sequence-discrete population investigated. \"\n$part1 += $lineBreak\n$part1 += \"I am also currently working on a manuscript concerning the population dynamics of the ubiquitous and abundant freshwater bacterium acI. \"\n$part1 += \"We have 14 single cell genomes(SAGs) from this clade, representing 3 lakes and X tribes. \"\n$part1 += \"Since these genomes come from a single cell and not an assembly of many cells, like the GFMs, we shredded the genomes and mapped them against each other in order to understand how their relationships would hold up in our metagenomic time series. \"\n$part1 += \"We found that SAGs from different tribes do not map well to each other. \"\n$part1 += \"When we mapped reads from the metagenomic time series of Mendota, we found that only SAGs colleted from the same lake were representive of the populations in the lake. \"\n$part1 += \"We also found that members of the same tribe represent ecologically distinct populations in the lake. \"\n$part1 += $lineBreak\n\n$newPara.Range.Text = $part1\n\n# --- 2) Split the Oral Presentations talk list into two lines ---\n# Re-locate this paragraph now (after the insertion above), since paragraph\n# references are positional and the prior insert shifted everything after it.\n$oralPresentations = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t.StartsWith('\"Genome-wide Selective Sweeps') -and $oralPresentations -eq $null) {\n        $oralPresentations = $p\n    }\n}\nif ($oralPresentations -eq $null) {\n    throw \"Could not find the Oral Presentations talks paragraph\"\n}\n\n$talk1 = '\"Genome-wide Selective Sweeps in Natural Bacterial Populations Revealed by Time-series Metagenomics.\" 15th International Symposium on Microbial Ecology. August 24, 2014. Seoul, South Korea'\n$talk2 = '\"Genome-wide and Gene-specific Selective Sweeps in Freshwater Bacterial Populations Revealed Using Metagenomics.\" 14 Symposium Society for Aquatic Microbial Ecology. August 2015. Uppsala, Sweden'\n$oralPresentations.Range.Text = $talk1 + $lineBreak + $talk2\n\nWrite-Output \"done\"\n"}
